$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (days since 1899-12-30).
# All data rows (2 through 18) need their serial value incremented by 1
# (45204 -> 45205, i.e. 2023-10-05 -> 2023-10-06).
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
